# Updated symbol list on Tue Jan  3 06:38:02 UTC 2023 with GitHub Actions
#
# Refreshes the crypto price/volume snapshot (Price column D, Volume(1h)
# column E) to the latest values pulled by the scraper, preserving the
# existing "plain text" cell formatting (no numeric auto-conversion).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'246.24"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'1.06%"
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'29.50"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'-1.42%"
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'5.153"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'0.06%"
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'0.05797"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'2.21%"
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'6.648"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'1.67%"
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'3.193"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'5.33%"
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.8516"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'0.47%"
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.8654"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'0.49%"
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.1375"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'1.81%"
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.07063"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'2.00%"
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.03262"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'13.25%"
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'0.09367"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'-0.14%"
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'0.001526"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'0.82%"
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'0.0005969"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'-0.13%"
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'0.006092"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'-0.31%"
$ws.Range("E16").Style = "Normal"
$ws.Range("E17").Value = "'-0.59%"
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'2.223"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'-0.48%"
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'0.3199"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'1.59%"
$ws.Range("E19").Style = "Normal"
$ws.Range("D21").Value = "'0.1281"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'-1.66%"
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'3.133"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'-13.73%"
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'0.04140"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'-0.75%"
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'0.1378"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'0.38%"
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'0.001227"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'1.36%"
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'0.004146"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'-6.69%"
$ws.Range("E26").Style = "Normal"
$ws.Range("E27").Value = "'2.44%"
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'0.0001442"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'3.73%"
$ws.Range("E28").Style = "Normal"
$ws.Range("D40").Value = "'0.03745"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'-0.03%"
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.005792"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'-0.09%"
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.1073"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'1.41%"
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.002447"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'19.09%"
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.009138"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'-1.62%"
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.00005270"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'3.13%"
$ws.Range("E45").Style = "Normal"
$ws.Range("E46").Value = "'-0.15%"
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'0.05789"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'-42.08%"
$ws.Range("E47").Style = "Normal"
$ws.Range("E48").Value = "'-21.66%"
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'0.00002096"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'-0.15%"
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'0.0001996"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'-0.15%"
$ws.Range("E50").Style = "Normal"
